$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Values are set with a leading apostrophe to force Excel to treat them as
# literal text (preventing values like "19.60" or "1.007" from being
# reinterpreted as numbers and losing trailing zeros / dot-grouping).
$ws.Range("D2").Value = "'26.888.36"
$ws.Range("D3").Value = "'1.826.88"
$ws.Range("D4").Value = "'1.007"
$ws.Range("D5").Value = "'310.68"
$ws.Range("D7").Value = "'0.4572"
$ws.Range("D8").Value = "'0.3686"
$ws.Range("D9").Value = "'0.07152"
$ws.Range("D10").Value = "'0.8744"
$ws.Range("D11").Value = "'0.07771"
$ws.Range("D12").Value = "'19.60"
$ws.Range("D13").Value = "'1.846.65"
$ws.Range("D14").Value = "'5.317"
$ws.Range("D15").Value = "'6.384"
$ws.Range("D16").Value = "'86.86"
$ws.Range("D17").Value = "'1.008"
$ws.Range("D18").Value = "'0.000008719"
$ws.Range("D20").Value = "'26.909.32"
$ws.Range("D22").Value = "'5.000"
$ws.Range("D23").Value = "'2.052.28"
$ws.Range("D24").Value = "'10.43"
$ws.Range("D25").Value = "'1.999"
$ws.Range("D26").Value = "'151.39"
$ws.Range("D28").Value = "'1.957"
$ws.Range("D29").Value = "'113.72"
$ws.Range("D30").Value = "'4.896"
$ws.Range("D31").Value = "'0.08798"
$ws.Range("D32").Value = "'3.045"
$ws.Range("D33").Value = "'0.7485"
$ws.Range("D34").Value = "'4.485"
$ws.Range("D35").Value = "'1.133"
$ws.Range("D36").Value = "'2.531"
$ws.Range("D38").Value = "'0.01943"
$ws.Range("D39").Value = "'0.05136"
$ws.Range("D40").Value = "'2.915"
$ws.Range("D41").Value = "'6.933"
$ws.Range("D42").Value = "'0.4964"
$ws.Range("D43").Value = "'0.1596"
$ws.Range("D44").Value = "'8.319"
$ws.Range("D45").Value = "'0.4692"
$ws.Range("D47").Value = "'10.16"
$ws.Range("D48").Value = "'101.37"
$ws.Range("D50").Value = "'0.06100"

# Strip the "quote prefix" text-format styling that Excel applies when a
# leading apostrophe is used, restoring cells to their original (unstyled)
# appearance while keeping the values as text.
$ws.Range("D2:D51").ClearFormats()

# --- Column E (Volume 1h %) updates ---
# These already contain non-numeric characters (surrounding spaces, '%')
# so Excel keeps them as plain text without any special handling needed.
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("E16").Value = "  -5.45%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  +4.27%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  -5.32%  "
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("E51").Value = "  -2.03%  "
